$d = $word.ActiveDocument

# The target text "constel·lació, Constel·lació" should become " Constel·lació"
# (dropping "constel·lació," and keeping a leading space, yielding a double space).
# We locate the substring with Find (no ReplaceWith) and overwrite the matched
# Range.Text directly; this avoids Word's Find/Replace "smart quotes" AutoFormat
# from mangling the straight apostrophes (') elsewhere in the same sentence.

$searchText = "constel·lació, Constel·lació"
$newText = " Constel·lació"

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

while ($rng.Find.Found) {
    $rng.Text = $newText
    $rng.Collapse(0)
    $rng.End = $d.Content.End
    $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}
